# Update TPM-derived NATMI metrics on the active sheet.
# Columns: G,H,I,J = Ligand avg/total expr + derived specificity (avg/total)
#          M,N,O,P = Receptor avg/total expr + derived specificity (avg/total)
#          Q,R,S,T = Edge avg/total expr weight + derived specificity (avg/total)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 17.16209533333333;  "H2" = 51.486286;         "I2" = 0.2459970657298922; "J2" = 0.2459970657298922
    "M2" = 174.1282373333333;  "N2" = 522.384712;        "O2" = 0.985625830323027;  "P2" = 0.985625830323027
    "Q2" = 2988.405409339959;  "R2" = 26895.64868405963; "S2" = 0.2424610621670533; "T2" = 0.2424610621670532

    "G3" = 17.16209533333333;  "H3" = 51.486286;         "I3" = 0.2459970657298922; "J3" = 0.2459970657298922
    "O3" = 0.003686901313133159; "P3" = 0.003686901313133159
    "Q3" = 11.17863948863711;  "R3" = 100.607755397734;  "S3" = 0.0009069669046664436; "T3" = 0.0009069669046664435

    "G4" = 17.16209533333333;  "H4" = 51.486286;         "I4" = 0.2459970657298922; "J4" = 0.2459970657298922
    "M4" = 1.888095;           "N4" = 5.664285;          "O4" = 0.01068726836383999; "P4" = 0.01068726836383999
    "Q4" = 32.40366638839;     "R4" = 291.63299749551;   "S4" = 0.002629036658172545; "T4" = 0.002629036658172545

    "I5" = 0.2105756965403629; "J5" = 0.2105756965403628
    "M5" = 174.1282373333333;  "N5" = 522.384712;        "O5" = 0.985625830323027;  "P5" = 0.985625830323027
    "Q5" = 2558.101856823419;  "R5" = 23022.91671141078; "S5" = 0.2075488457484449; "T5" = 0.2075488457484449

    "I6" = 0.2105756965403629; "J6" = 0.2105756965403628
    "O6" = 0.003686901313133159; "P6" = 0.003686901313133159
    "Q6" = 9.569015751098554;  "R6" = 86.121141759887;   "S6" = 0.0007763718120885934; "T6" = 0.0007763718120885932

    "I7" = 0.2105756965403629; "J7" = 0.2105756965403628
    "M7" = 1.888095;           "N7" = 5.664285;          "O7" = 0.01068726836383999; "P7" = 0.01068726836383999
    "Q7" = 27.737829310895;    "R7" = 249.640463798055;  "S7" = 0.002250478979829391; "T7" = 0.002250478979829391

    "G8" = 37.91244433333333;  "H8" = 113.737333;        "I8" = 0.543427237729745;  "J8" = 0.543427237729745
    "M8" = 174.1282373333333;  "N8" = 522.384712;        "O8" = 0.985625830323027;  "P8" = 0.985625830323027
    "Q8" = 6601.627104761455;  "R8" = 59414.6439428531;  "S8" = 0.5356159224075289; "T8" = 0.5356159224075289

    "G9" = 37.91244433333333;  "H9" = 113.737333;        "I9" = 0.543427237729745;  "J9" = 0.543427237729745
    "O9" = 0.003686901313133159; "P9" = 0.003686901313133159
    "Q9" = 24.69451072866411;  "R9" = 222.250596557977;  "S9" = 0.002003562596378122; "T9" = 0.002003562596378122

    "G10" = 37.91244433333333; "H10" = 113.737333;       "I10" = 0.543427237729745; "J10" = 0.543427237729745
    "M10" = 1.888095;          "N10" = 5.664285;         "O10" = 0.01068726836383999; "P10" = 0.01068726836383999
    "Q10" = 71.582296583545;   "R10" = 644.240669251905; "S10" = 0.005807752725838059; "T10" = 0.005807752725838059
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
